# Helper: force the engine to keep a sequence of adjacent character runs as
# *separate* <w:r> elements (all sharing identical, empty <w:rPr/>) instead of
# silently re-merging them into one run. Toggling a formatting property on and
# back off at each boundary is enough to make the run-split "stick" without
# leaving any residual formatting behind.
function Split-RunsAt($doc, $rangeStart, $boundaries) {
    foreach ($b in $boundaries) {
        $sub = $doc.Range($rangeStart, $b)
        $sub.Bold = 1
        $sub2 = $doc.Range($rangeStart, $b)
        $sub2.Bold = 0
    }
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update the "Errata- Last modified on ..." date on the very first
#    paragraph: 1/23/2019 -> 9/12/2019, re-splitting into the run layout
#    used by the edited document.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$textOnly = $d.Range($p1.Start, $p1.End - 1)
$textOnly.Delete()

$p1 = $d.Paragraphs(1).Range
$p1.InsertAfter("Errata- Last modified on 9/12/2019")

$s = $p1.Start
$e = $p1.End
# "Errata- Last modified on " (25) | "9" (1) | "/" (1) | "12" (2) | "/2019" (5)
Split-RunsAt $d $s @($s + 25, $s + 26, $s + 27, $s + 29, $e)

# ---------------------------------------------------------------------------
# 2) Append the new "Page 159" errata entry after the existing
#    "To: but you do not have permission to write to it." paragraph. Four
#    new paragraphs are added: a blank separator, the "Page 159 ..." line,
#    the "From: ..." line (2 runs) and the "To: ..." line (6 runs).
# ---------------------------------------------------------------------------
$writePara = $d.Paragraphs(9).Range
$writePara.InsertParagraphAfter()

# Paragraph 10: blank separator - left empty, matches the diff.
$blankPara = $d.Paragraphs(10).Range
$blankPara.InsertParagraphAfter()

# Paragraph 11: "Page 159 ..." (single run).
$pagePara = $d.Paragraphs(11).Range
$pagePara.InsertAfter("Page 159 – Thanks to Tom Parsons for this correction!")
$pagePara.InsertParagraphAfter()

# Paragraph 12: "From: ..." (2 runs).
$fromPara = $d.Paragraphs(12).Range
$fromPara.InsertAfter("From: The following command recursively searches the root directory for all the files that have inode number 36700164")
$fromPara.InsertParagraphAfter()
$fs = $fromPara.Start
$fe = $fromPara.End
# "From: " (6)
Split-RunsAt $d $fs @($fs + 6, $fe)

# Paragraph 13: "To: ..." (6 runs).
$toPara = $d.Paragraphs(13).Range
$toPara.InsertAfter("To: The following command recursively searches the current and /usr directories for all the files that have inode number 36700164")
$ts = $toPara.Start
$te = $toPara.End
# "To: " (4) | "The following command recursively searches the " (49->51) | "current and /usr " (17->68) | "director" (8->76) | "ies" (3->79) | " for all the files that have inode number 36700164" (50->129)
Split-RunsAt $d $ts @($ts + 4, $ts + 51, $ts + 68, $ts + 76, $ts + 79, $te)
